# Generate Report for Handoff
#
# - The pending handoff for 48f06b51-765c-4487-b49f-b8a7fa7a4f33.md completed:
#   it is renamed/re-identified as 8664dfb4-cf42-40a2-ab26-496d7a14a04a.md, with
#   a fresh handoff package (new content hash, new handoff timestamps).
# - A new source file c7ecd104-466a-40be-ab38-08c768c7283c.md was picked up but
#   its handoff transform failed, so it is listed with status
#   "Handoff transform failed" / dependency "Ignored".
# - The ".localization-config" row (not localized) shifts down one row on every
#   sheet to make room for the new file.

$wb = $excel.ActiveWorkbook

$hyperLinkColor = 15570276  # BGR for #6495ED - matches the workbook's existing HyperLink font color

function Set-HyperlinkStyle($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperLinkColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Drop every existing hyperlink on the sheet so we can rebuild them (and their
# relationship ids) cleanly in the new order.
$ws1.Range("A1:C3").Hyperlinks.Delete()

$ws1.Range("A2").Value = "8664dfb4-cf42-40a2-ab26-496d7a14a04a.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("A3").Value = "c7ecd104-466a-40be-ab38-08c768c7283c.md"
$ws1.Range("B3").Value = "Handoff transform failed"
$ws1.Range("C3").Value = "Handoff transform failed"

$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ee0154dffe9382ba9b29d6bc8b4c06e0a00758fe/e2e/8664dfb4-cf42-40a2-ab26-496d7a14a04a.md", "", "", "8664dfb4-cf42-40a2-ab26-496d7a14a04a.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ee0154dffe9382ba9b29d6bc8b4c06e0a00758fe/e2e/c7ecd104-466a-40be-ab38-08c768c7283c.md", "", "", "c7ecd104-466a-40be-ab38-08c768c7283c.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ee0154dffe9382ba9b29d6bc8b4c06e0a00758fe/.localization-config", "", "", ".localization-config") | Out-Null

Set-HyperlinkStyle $ws1.Range("A2")
Set-HyperlinkStyle $ws1.Range("A3")
Set-HyperlinkStyle $ws1.Range("A4")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A1:I3").Hyperlinks.Delete()

$ws2.Range("A2").Value = "8664dfb4-cf42-40a2-ab26-496d7a14a04a.md"
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = "8664dfb4-cf42-40a2-ab26-496d7a14a04a.2dbe59cb6d0dadcfef2c5f36b85ad0e652679c64.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-26 09:27:00"
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "c7ecd104-466a-40be-ab38-08c768c7283c.md"
$ws2.Range("B3").Value = "Handoff transform failed"
$ws2.Range("D3").Value = "0001-01-01 00:00:00"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Ignored"

$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Ignored"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ee0154dffe9382ba9b29d6bc8b4c06e0a00758fe/e2e/8664dfb4-cf42-40a2-ab26-496d7a14a04a.md", "", "", "8664dfb4-cf42-40a2-ab26-496d7a14a04a.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8e350483b3d023677d0b59cf0aa6a836d101b363/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/8664dfb4-cf42-40a2-ab26-496d7a14a04a.2dbe59cb6d0dadcfef2c5f36b85ad0e652679c64.zh-cn.xlf", "", "", "8664dfb4-cf42-40a2-ab26-496d7a14a04a.2dbe59cb6d0dadcfef2c5f36b85ad0e652679c64.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ee0154dffe9382ba9b29d6bc8b4c06e0a00758fe/e2e/c7ecd104-466a-40be-ab38-08c768c7283c.md", "", "", "c7ecd104-466a-40be-ab38-08c768c7283c.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ee0154dffe9382ba9b29d6bc8b4c06e0a00758fe/.localization-config", "", "", ".localization-config") | Out-Null

Set-HyperlinkStyle $ws2.Range("A2")
Set-HyperlinkStyle $ws2.Range("C2")
Set-HyperlinkStyle $ws2.Range("A3")
Set-HyperlinkStyle $ws2.Range("A4")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A1:I3").Hyperlinks.Delete()

$ws3.Range("A2").Value = "8664dfb4-cf42-40a2-ab26-496d7a14a04a.md"
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = "8664dfb4-cf42-40a2-ab26-496d7a14a04a.2dbe59cb6d0dadcfef2c5f36b85ad0e652679c64.de-de.xlf"
$ws3.Range("D2").Value = "2016-01-26 09:27:14"
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "c7ecd104-466a-40be-ab38-08c768c7283c.md"
$ws3.Range("B3").Value = "Handoff transform failed"
$ws3.Range("D3").Value = "0001-01-01 00:00:00"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Ignored"

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Ignored"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/ee0154dffe9382ba9b29d6bc8b4c06e0a00758fe/e2e/8664dfb4-cf42-40a2-ab26-496d7a14a04a.md", "", "", "8664dfb4-cf42-40a2-ab26-496d7a14a04a.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cdc0967369f7f30ed7fc71a4f59ac0d9138cee25/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/8664dfb4-cf42-40a2-ab26-496d7a14a04a.2dbe59cb6d0dadcfef2c5f36b85ad0e652679c64.de-de.xlf", "", "", "8664dfb4-cf42-40a2-ab26-496d7a14a04a.2dbe59cb6d0dadcfef2c5f36b85ad0e652679c64.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/ee0154dffe9382ba9b29d6bc8b4c06e0a00758fe/e2e/c7ecd104-466a-40be-ab38-08c768c7283c.md", "", "", "c7ecd104-466a-40be-ab38-08c768c7283c.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/ee0154dffe9382ba9b29d6bc8b4c06e0a00758fe/.localization-config", "", "", ".localization-config") | Out-Null

Set-HyperlinkStyle $ws3.Range("A2")
Set-HyperlinkStyle $ws3.Range("C2")
Set-HyperlinkStyle $ws3.Range("A3")
Set-HyperlinkStyle $ws3.Range("A4")

Write-Output "Report generated for handoff"
